$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Add new row 18: npv_with_annuity parameter
$ws.Range("A18").Value = "npv_with_annuity"
$ws.Range("B18").Value = $true
$ws.Range("C18").Value = "If this is true, the npv is calculated with the annuities, not with the restpayment"

# Update the active selection on this sheet to C3 (as reflected in the saved view state)
$ws.Activate()
$ws.Range("C3").Select()
